# Daily attendance processing - 2025-10-10 21:59:23
# Applies updates to the "Session Analysis Results" sheet:
#  - reorders several "Recorded By" email lists
#  - updates Missing/Pending session counters
#  - updates the BIOCHEMISTRY LAB/CBL session row (date, status, style)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# --- Reorder "Recorded By" email lists (G column) ---
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G12").Value = "mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
$ws.Range("G25").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G26").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G34").Value = "mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
$ws.Range("G41").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# --- Update class statistics counters ---
$ws.Range("L7").Value = 4    # Missing Sessions
$ws.Range("L8").Value = 33   # Pending Sessions

# --- Update Year 3 / C2 group statistics row ---
$ws.Range("P16").Value = 2   # Missing
$ws.Range("Q16").Value = 16  # Pending

# --- Update row 33 (Year 3, C2, BIOCHEMISTRY LAB/CBL, session 1) ---
# Date corrected from 07/11/2025 to 07/10/2025 and status moved from
# Pending to Not Recorded. Force the date cell to be stored as text (not
# auto-converted to a date serial) before writing the new value.
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "07/10/2025"
$ws.Range("I33").Value = "Not Recorded"

# Re-style the whole row to match the other "Not Recorded" rows (e.g. row 24)
$ws.Range("A24:I24").Copy()
$ws.Range("A33:I33").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
